$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Workbook / worksheet metadata updates ---
$ws.Name = "IClientBalance-20240820-092050-"

# --- Column G (report date): 2024-08-19 (45523) -> 2024-08-20 (45524) for rows 2..274 ---
for ($r = 2; $r -le 274; $r++) {
    $ws.Cells.Item($r, 7).Value = 45524
}

# --- Row-specific balance updates (columns D, E, H) ---
# Row 5
$ws.Cells.Item(5, 4).Value = 18107.37
$ws.Cells.Item(5, 8).Value = 18850.7
# Row 8
$ws.Cells.Item(8, 4).Value = 6853.85
$ws.Cells.Item(8, 8).Value = 9877.61
# Row 15
$ws.Cells.Item(15, 4).Value = 23658.2
$ws.Cells.Item(15, 5).Value = 10135.48
$ws.Cells.Item(15, 8).Value = 33793.68
# Row 17
$ws.Cells.Item(17, 4).Value = 8007.03
$ws.Cells.Item(17, 5).Value = 3442.08
$ws.Cells.Item(17, 8).Value = 11449.11
# Row 43
$ws.Cells.Item(43, 4).Value = 12728.81
$ws.Cells.Item(43, 5).Value = 5016.22
$ws.Cells.Item(43, 8).Value = 17745.03
# Row 49
$ws.Cells.Item(49, 4).Value = 5127.47
$ws.Cells.Item(49, 8).Value = 6915.33
# Row 60
$ws.Cells.Item(60, 4).Value = 16165.02
$ws.Cells.Item(60, 8).Value = 16926.080000000002
# Row 99
$ws.Cells.Item(99, 4).Value = 13951.72
$ws.Cells.Item(99, 5).Value = 5162.29
$ws.Cells.Item(99, 8).Value = 19114.009999999998
# Row 102
$ws.Cells.Item(102, 5).Value = 10987.23
$ws.Cells.Item(102, 8).Value = 10987.23
# Row 104
$ws.Cells.Item(104, 4).Value = 44993.99
$ws.Cells.Item(104, 5).Value = 15167.95
$ws.Cells.Item(104, 8).Value = 60161.94
# Row 105
$ws.Cells.Item(105, 5).Value = 753.72
$ws.Cells.Item(105, 8).Value = 753.72
# Row 108
$ws.Cells.Item(108, 4).Value = 46020.25
$ws.Cells.Item(108, 5).Value = 14726.6
$ws.Cells.Item(108, 8).Value = 60746.85
# Row 129
$ws.Cells.Item(129, 4).Value = 2052.7800000000002
$ws.Cells.Item(129, 8).Value = 2177.36
# Row 132
$ws.Cells.Item(132, 4).Value = 5709.53
$ws.Cells.Item(132, 5).Value = 2931.08
$ws.Cells.Item(132, 8).Value = 8640.61
# Row 143
$ws.Cells.Item(143, 4).Value = 48697.11
$ws.Cells.Item(143, 5).Value = 17403.97
$ws.Cells.Item(143, 8).Value = 66101.08
# Row 158
$ws.Cells.Item(158, 4).Value = 1070.3399999999999
$ws.Cells.Item(158, 5).Value = 1548
$ws.Cells.Item(158, 8).Value = 2618.34
# Row 173
$ws.Cells.Item(173, 4).Value = 27537.75
$ws.Cells.Item(173, 5).Value = 9232.42
$ws.Cells.Item(173, 8).Value = 36770.17
# Row 231
$ws.Cells.Item(231, 5).Value = 775.95
$ws.Cells.Item(231, 8).Value = 775.95
# Row 232
$ws.Cells.Item(232, 5).Value = 46018.17
$ws.Cells.Item(232, 8).Value = 46018.17
# Row 235
$ws.Cells.Item(235, 4).Value = 11209.39
$ws.Cells.Item(235, 5).Value = 3802.88
$ws.Cells.Item(235, 8).Value = 15012.27
# Row 264
$ws.Cells.Item(264, 4).Value = 34212.639999999999
$ws.Cells.Item(264, 5).Value = 11543.76
$ws.Cells.Item(264, 8).Value = 45756.4
# Row 265
$ws.Cells.Item(265, 4).Value = 18531.02
$ws.Cells.Item(265, 5).Value = 6695.35
$ws.Cells.Item(265, 8).Value = 25226.37
# Row 270
$ws.Cells.Item(270, 4).Value = 13941.48
$ws.Cells.Item(270, 8).Value = 14816.2
# Row 271
$ws.Cells.Item(271, 4).Value = 18308.72
$ws.Cells.Item(271, 8).Value = 19177.259999999998
# Row 273
$ws.Cells.Item(273, 4).Value = 11965.96
$ws.Cells.Item(273, 5).Value = 4677.6899999999996
$ws.Cells.Item(273, 8).Value = 16643.650000000001
